$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-07-12 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-13 Thursday", 2)

# Update the 20x5 table of arithmetic expressions, in row-major order
$newValues = @(
    "64+34=",
    "64-7=",
    "38+1=",
    "3-0=",
    "81-42=",
    "4+17=",
    "74+24=",
    "14+8=",
    "71-51=",
    "2+26=",
    "23+70=",
    "87-31=",
    "43-8=",
    "45+17=",
    "46+3=",
    "80-66=",
    "19+26=",
    "67+16=",
    "50-45=",
    "14+13=",
    "35-33=",
    "25+49=",
    "70-40=",
    "1+15=",
    "29-21=",
    "62-3=",
    "78-77=",
    "45+34=",
    "16+31=",
    "27+19=",
    "62-24=",
    "42+32=",
    "22+35=",
    "14-4=",
    "1+69=",
    "57+11=",
    "70-63=",
    "5+6=",
    "64-8=",
    "50+5=",
    "27+60=",
    "43-22=",
    "76+15=",
    "87-67=",
    "21-20=",
    "43+48=",
    "44+13=",
    "93-8=",
    "74+13=",
    "61+34=",
    "37-14=",
    "0+34=",
    "61+13=",
    "17+23=",
    "10-7=",
    "37-5=",
    "31-16=",
    "11+42=",
    "27+53=",
    "2+9=",
    "60-26=",
    "10+49=",
    "78-74=",
    "82-24=",
    "84+8=",
    "64-3=",
    "64-40=",
    "76+18=",
    "45-5=",
    "31-13=",
    "90-63=",
    "27+4=",
    "50-6=",
    "32+12=",
    "11+8=",
    "92-70=",
    "40-10=",
    "46-29=",
    "52+31=",
    "54+25=",
    "76-25=",
    "57-25=",
    "65-31=",
    "67-11=",
    "2+42=",
    "67-62=",
    "61+28=",
    "67-1=",
    "37+48=",
    "0+2=",
    "84-28=",
    "19+38=",
    "85+8=",
    "59-58=",
    "32-4=",
    "98-45=",
    "4+31=",
    "70-25=",
    "68-9=",
    "54-3="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Done. Updated" $idx "cells."